# Results from R script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 103: correct the timestamp (A103) -- data revision, same date, different time.
$ws.Range("A103").Value = 45470.2916666667

# Row 104: append the next day's OHLCV record.
# Seed A104 from A103 first so it inherits the same date/time cell style (s="1")
# instead of minting a brand-new style entry.
$ws.Range("A103").Copy($ws.Range("A104"))
$ws.Range("A104").Value = 45471.3652430556

$ws.Range("B104").Value = 23769
$ws.Range("C104").Value = 0.740000009536743
$ws.Range("D104").Value = 0.670000016689301
$ws.Range("E104").Value = 0.694999992847443
$ws.Range("F104").Value = 0.709999978542328

# G104 is text (looks numeric, stored as a shared string) -- seed it by copying an
# existing cell holding the identical string so the run reuses the shared-string
# entry and default style rather than creating a new one.
$ws.Range("G95").Copy($ws.Range("G104"))

# H104 ticker -- plain text, reuses the existing "BWZ.MI" shared string.
$ws.Range("H104").Value = "BWZ.MI"
